$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numericCells = @("D5", "D6", "D9", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D27", "D29", "D33", "D38", "D40", "D41", "D43", "D44", "D45", "D46", "D48", "D50", "D51")
foreach ($addr in $numericCells) { $ws.Range($addr).NumberFormat = '@' }

$ws.Range("D2").Value = '66.109.35'
$ws.Range("E2").Value = '  +0.46%  '
$ws.Range("D3").Value = '3.558.75'
$ws.Range("E3").Value = '  +4.19%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '605.78'
$ws.Range("E5").Value = '  +1.58%  '
$ws.Range("D6").Value = '144.55'
$ws.Range("E6").Value = '  +1.90%  '
$ws.Range("D7").Value = '3.557.05'
$ws.Range("E7").Value = '  +4.19%  '
$ws.Range("D9").Value = '0.492'
$ws.Range("E9").Value = '  +4.84%  '
$ws.Range("E10").Value = '  +2.28%  '
$ws.Range("E11").Value = '  -1.42%  '
$ws.Range("E12").Value = '  +1.72%  '
$ws.Range("D13").Value = '4.160.75'
$ws.Range("E13").Value = '  +4.19%  '
$ws.Range("E14").Value = '  +3.79%  '
$ws.Range("E15").Value = '  +1.30%  '
$ws.Range("D16").Value = '3.559.20'
$ws.Range("E16").Value = '  +4.36%  '
$ws.Range("D17").Value = '66.184.67'
$ws.Range("E17").Value = '  +0.58%  '
$ws.Range("D18").Value = '0.116'
$ws.Range("E18").Value = '  -0.45%  '
$ws.Range("D19").Value = '11.30'
$ws.Range("E19").Value = '  +9.57%  '
$ws.Range("D20").Value = '6.18'
$ws.Range("E20").Value = '  +1.56%  '
$ws.Range("D21").Value = '14.82'
$ws.Range("E21").Value = '  +2.16%  '
$ws.Range("D22").Value = '428.89'
$ws.Range("E22").Value = '  +3.49%  '
$ws.Range("D23").Value = '0.611'
$ws.Range("E23").Value = '  +6.30%  '
$ws.Range("D24").Value = '79.22'
$ws.Range("E24").Value = '  +2.65%  '
$ws.Range("D25").Value = '3.698.47'
$ws.Range("E25").Value = '  +4.25%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").Value = '0.0000118'
$ws.Range("E27").Value = '  +7.89%  '
$ws.Range("E28").Value = '  +3.52%  '
$ws.Range("D29").Value = '7.95'
$ws.Range("E29").Value = '  +0.26%  '
$ws.Range("E30").Value = '  -2.38%  '
$ws.Range("E31").Value = '  +0.03%  '
$ws.Range("E32").Value = '  +0.79%  '
$ws.Range("D33").Value = '25.56'
$ws.Range("E33").Value = '  +4.00%  '
$ws.Range("D34").Value = '3.552.49'
$ws.Range("E34").Value = '  +4.11%  '
$ws.Range("E35").Value = '  -5.68%  '
$ws.Range("E36").Value = '  +0.07%  '
$ws.Range("E37").Value = '  +3.71%  '
$ws.Range("D38").Value = '7.86'
$ws.Range("E38").Value = '  +4.99%  '
$ws.Range("E39").Value = '  +2.28%  '
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("D41").Value = '175.36'
$ws.Range("E41").Value = '  +3.93%  '
$ws.Range("E42").Value = '  -0.35%  '
$ws.Range("D43").Value = '5.20'
$ws.Range("E43").Value = '  +3.39%  '
$ws.Range("D44").Value = '0.893'
$ws.Range("E44").Value = '  +2.50%  '
$ws.Range("D45").Value = '1.94'
$ws.Range("E45").Value = '  +2.25%  '
$ws.Range("D46").Value = '46.07'
$ws.Range("E46").Value = '  +1.54%  '
$ws.Range("E47").Value = '  +0.36%  '
$ws.Range("D48").Value = '25.67'
$ws.Range("E48").Value = '  -2.40%  '
$ws.Range("E49").Value = '  +16.13%  '
$ws.Range("B50").Value = 'dogwifhat'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D50").Value = '2.34'
$ws.Range("E50").Value = '  +3.00%  '
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").Value = '7.11'
$ws.Range("E51").Value = '  +0.95%  '

Write-Host "Applied cryptos list update"
